# Auto-generated Excel COM-interop script to apply the Leviathan_Profits.xlsx data refresh
# (static market-price / profit values recalculated by the scheduled runner).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 7686.25
$ws.Range("I34").Value = 7686.25
$ws.Range("K34").Value = 7686.25
$ws.Range("M34").Value = -7483.25
$ws.Range("H36").Value = 7686.25
$ws.Range("I36").Value = 7686.25
$ws.Range("K36").Value = 7686.25
$ws.Range("M36").Value = -6971.25
$ws.Range("H100").Value = 2783.0833
$ws.Range("I100").Value = 2232.5557
$ws.Range("K100").Value = 2232.5557
$ws.Range("M100").Value = -1691.5557
$ws.Range("H106").Value = 6408.0415
$ws.Range("J106").Value = 17538.428
$ws.Range("L106").Value = 17538.428
$ws.Range("N106").Value = -18800.428
$ws.Range("H112").Value = 1792.3334
$ws.Range("J112").Value = 1854.7368
$ws.Range("L112").Value = 5564.2104
$ws.Range("N112").Value = -7780.2104
$ws.Range("H137").Value = 37087.75
$ws.Range("I137").Value = 1427.8334
$ws.Range("J137").Value = 251047.25
$ws.Range("K137").Value = 4283.5002
$ws.Range("L137").Value = 753141.75
$ws.Range("M137").Value = -1733.5002
$ws.Range("N137").Value = -758241.75
$ws.Range("H138").Value = 2116
$ws.Range("J138").Value = 2997.75
$ws.Range("L138").Value = 8993.25
$ws.Range("N138").Value = -19273.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5999.8335
$ws.Range("I61").Value = 5999.8335
$ws.Range("K61").Value = 5999.8335
$ws.Range("M61").Value = -5787.8335
$ws.Range("H63").Value = 2416.4443
$ws.Range("J63").Value = 2979.8
$ws.Range("L63").Value = 2979.8
$ws.Range("N63").Value = -4351.8
$ws.Range("H66").Value = 2416.4443
$ws.Range("J66").Value = 2979.8
$ws.Range("L66").Value = 14899
$ws.Range("N66").Value = -21763
$ws.Range("H122").Value = 16736.8
$ws.Range("I122").Value = 17646.572
$ws.Range("K122").Value = 52939.716
$ws.Range("M122").Value = -50489.716
$ws.Range("H136").Value = 5999.8335
$ws.Range("I136").Value = 5999.8335
$ws.Range("K136").Value = 17999.5005
$ws.Range("M136").Value = -15449.5005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15628.5
$ws.Range("H85").Value = 15628.5
$ws.Range("H99").Value = 1336.1
$ws.Range("I99").Value = 1336.1
$ws.Range("K99").Value = 1336.1
$ws.Range("M99").Value = 161.9000000000001
$ws.Range("H105").Value = 3707881.5
$ws.Range("J105").Value = 2834.5557
$ws.Range("L105").Value = 2834.5557
$ws.Range("N105").Value = -6328.5557
$ws.Range("H107").Value = 12416.9375
$ws.Range("I107").Value = 3667.1538
$ws.Range("K107").Value = 3667.1538
$ws.Range("M107").Value = -1747.1538
$ws.Range("H128").Value = 21250
$ws.Range("I128").Value = 21250
$ws.Range("K128").Value = 63750
$ws.Range("M128").Value = -61260
$ws.Range("H134").Value = 3034.3684
$ws.Range("I134").Value = 2610.2
$ws.Range("K134").Value = 7830.599999999999
$ws.Range("M134").Value = -5295.599999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1590.5
$ws.Range("I16").Value = 1457.5
$ws.Range("K16").Value = 1457.5
$ws.Range("M16").Value = -1170.5
$ws.Range("H105").Value = 6607.84
$ws.Range("I105").Value = 5206.9287
$ws.Range("K105").Value = 5206.9287
$ws.Range("M105").Value = -3459.9287
$ws.Range("H113").Value = 1590.5
$ws.Range("I113").Value = 1457.5
$ws.Range("K113").Value = 1457.5
$ws.Range("M113").Value = 712.5
$ws.Range("H132").Value = 2380.7354
$ws.Range("I132").Value = 2239.0688
$ws.Range("K132").Value = 6717.2064
$ws.Range("M132").Value = -4187.2064
$ws.Range("H141").Value = 86170.875
$ws.Range("I141").Value = 83333.336
$ws.Range("J141").Value = 87873.39999999999
$ws.Range("K141").Value = 83333.336
$ws.Range("L141").Value = 87873.39999999999
$ws.Range("M141").Value = -78153.336
$ws.Range("N141").Value = -98233.39999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 725
$ws.Range("J2").Value = 725
$ws.Range("L2").Value = 4350
$ws.Range("N2").Value = -4576
$ws.Range("H37").Value = 1000000000
$ws.Range("J37").Value = 1000000000
$ws.Range("L37").Value = 3000000000
$ws.Range("N37").Value = -3000000224
$ws.Range("H38").Value = 121.125
$ws.Range("I38").Value = 117.25
$ws.Range("K38").Value = 351.75
$ws.Range("M38").Value = -4.75
$ws.Range("H92").Value = 1657.5555
$ws.Range("I92").Value = 1200
$ws.Range("J92").Value = 2023.6
$ws.Range("K92").Value = 3600
$ws.Range("L92").Value = 6070.799999999999
$ws.Range("M92").Value = -2352
$ws.Range("N92").Value = -8566.799999999999
$ws.Range("H123").Value = 6865.6665
$ws.Range("I123").Value = 2798.5
$ws.Range("J123").Value = 15000
$ws.Range("K123").Value = 8395.5
$ws.Range("L123").Value = 45000
$ws.Range("M123").Value = -5945.5
$ws.Range("N123").Value = -49900
$ws.Range("H125").Value = 9974.5
$ws.Range("I125").Value = 9974.5
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 29923.5
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -25003.5
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 8996.333000000001
$ws.Range("I126").Value = 8994.5
$ws.Range("K126").Value = 26983.5
$ws.Range("M126").Value = -22043.5
$ws.Range("H132").Value = 1661.15
$ws.Range("J132").Value = 1947.3572
$ws.Range("L132").Value = 17526.2148
$ws.Range("N132").Value = -22586.2148

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10053616
$ws.Range("I11").Value = 18122428
$ws.Range("J11").Value = 639999.8
$ws.Range("K11").Value = 18122428
$ws.Range("L11").Value = 639999.8
$ws.Range("M11").Value = -18122289
$ws.Range("N11").Value = -640277.8
$ws.Range("H126").Value = 2602.5
$ws.Range("I126").Value = 2470.3333
$ws.Range("K126").Value = 7410.999899999999
$ws.Range("M126").Value = -4940.999899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5375.143
$ws.Range("I7").Value = 4521.421
$ws.Range("J7").Value = 7177.4443
$ws.Range("K7").Value = 4521.421
$ws.Range("L7").Value = 7177.4443
$ws.Range("M7").Value = -4409.421
$ws.Range("N7").Value = -7401.4443
$ws.Range("H16").Value = 1162.1333
$ws.Range("I16").Value = 1162.1333
$ws.Range("K16").Value = 1162.1333
$ws.Range("M16").Value = -992.1333
$ws.Range("H22").Value = 2178.75
$ws.Range("I22").Value = 1956.375
$ws.Range("K22").Value = 1956.375
$ws.Range("M22").Value = -1661.375
$ws.Range("H27").Value = 2178.75
$ws.Range("I27").Value = 1956.375
$ws.Range("K27").Value = 1956.375
$ws.Range("M27").Value = -1849.375
$ws.Range("H40").Value = 4399.5356
$ws.Range("I40").Value = 4324.3335
$ws.Range("K40").Value = 4324.3335
$ws.Range("M40").Value = -4188.3335
$ws.Range("H122").Value = 20084
$ws.Range("I122").Value = 35734.668
$ws.Range("K122").Value = 107204.004
$ws.Range("M122").Value = -104754.004
$ws.Range("H126").Value = 5375.143
$ws.Range("I126").Value = 4521.421
$ws.Range("J126").Value = 7177.4443
$ws.Range("K126").Value = 13564.263
$ws.Range("L126").Value = 21532.3329
$ws.Range("M126").Value = -11094.263
$ws.Range("N126").Value = -26472.3329
$ws.Range("H136").Value = 3123.1785
$ws.Range("I136").Value = 2413.4375
$ws.Range("J136").Value = 4069.5
$ws.Range("K136").Value = 7240.3125
$ws.Range("L136").Value = 12208.5
$ws.Range("M136").Value = -4690.3125
$ws.Range("N136").Value = -17308.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3200
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3200
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3200
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -3426
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H52").Value = 23873.25
$ws.Range("I52").Value = 18999
$ws.Range("J52").Value = 28747.5
$ws.Range("K52").Value = 18999
$ws.Range("L52").Value = 28747.5
$ws.Range("M52").Value = -18773
$ws.Range("N52").Value = -29199.5
$ws.Range("H132").Value = 5782.8335
$ws.Range("I132").Value = 5840.054
$ws.Range("K132").Value = 17520.162
$ws.Range("M132").Value = -14990.162
